$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.547008414586003
$ws.Range("C2").Value = 0.07789452405344832
$ws.Range("D2").Value = 0.07585075576982092
$ws.Range("E2").Value = 0.04605413356096832
$ws.Range("G2").Value = 0.002651283421994027
$ws.Range("K2").Value = 2.064225797974728
$ws.Range("L2").Value = 0.2704044349660819
$ws.Range("N2").Value = 5.3160120756271
$ws.Range("B3").Value = 2.482129439093057
$ws.Range("C3").Value = 0.06982563133261976
$ws.Range("D3").Value = 0.06921157424085322
$ws.Range("E3").Value = 0.04583233068188441
$ws.Range("G3").Value = 0.002657252036143642
$ws.Range("K3").Value = 1.996212624205413
$ws.Range("L3").Value = 0.2632962233936524
$ws.Range("N3").Value = 5.234119080627721
$ws.Range("B4").Value = 2.443829364750627
$ws.Range("C4").Value = 0.06491012946081298
$ws.Range("D4").Value = 0.06517865214294716
$ws.Range("E4").Value = 0.04571674914006607
$ws.Range("G4").Value = 0.002661106511512045
$ws.Range("K4").Value = 1.955752477497271
$ws.Range("L4").Value = 0.2590969603842126
$ws.Range("N4").Value = 5.184229373366747
$ws.Range("B5").Value = 2.428606437458257
$ws.Range("C5").Value = 0.06291653258986685
$ws.Range("D5").Value = 0.06354598158795
$ws.Range("E5").Value = 0.04567482925663846
$ws.Range("G5").Value = 0.002662725120790213
$ws.Range("K5").Value = 1.939589851928133
$ws.Range("L5").Value = 0.2574270923006452
$ws.Range("N5").Value = 5.163996299974428
$ws.Range("B6").Value = 2.426101870312834
$ws.Range("C6").Value = 0.06258606436009018
$ws.Range("D6").Value = 0.0632755237459719
$ws.Range("E6").Value = 0.04566818135796602
$ws.Range("G6").Value = 0.002662996786257241
$ws.Range("K6").Value = 1.936925652241683
$ws.Range("L6").Value = 0.2571523052672973
$ws.Range("N6").Value = 5.160642454335914
$ws.Range("B7").Value = 2.443622507573423
$ws.Range("C7").Value = 0.06488320494260336
$ws.Range("D7").Value = 0.06515658996578111
$ws.Range("E7").Value = 0.04571616281904411
$ws.Range("G7").Value = 0.002661128146284159
$ws.Range("K7").Value = 1.955533188184177
$ws.Range("L7").Value = 0.2590742727000901
$ws.Range("N7").Value = 5.183956110679361
$ws.Range("B8").Value = 2.524318405676524
$ws.Range("C8").Value = 0.0751041542181099
$ws.Range("D8").Value = 0.07355241597844042
$ws.Range("E8").Value = 0.04597337816994695
$ws.Range("G8").Value = 0.002653302127495127
$ws.Range("K8").Value = 2.040503914449744
$ws.Range("L8").Value = 0.2679191242568777
$ws.Range("N8").Value = 5.28769293944427
$ws.Range("B9").Value = 2.694831807086018
$ws.Range("C9").Value = 0.09546807421580183
$ws.Range("D9").Value = 0.09037113022728249
$ws.Range("E9").Value = 0.04664142096421919
$ws.Range("G9").Value = 0.002639452809570968
$ws.Range("K9").Value = 2.21753863801888
$ws.Range("L9").Value = 0.2865840669068831
$ws.Range("N9").Value = 5.494312885516422
$ws.Range("B10").Value = 2.827720222039716
$ws.Range("C10").Value = 0.1106439127221108
$ws.Range("D10").Value = 0.1029576105584056
$ws.Range("E10").Value = 0.04723232611098283
$ws.Range("G10").Value = 0.002630179568036741
$ws.Range("K10").Value = 2.354094252442167
$ws.Range("L10").Value = 0.3011170051718182
$ws.Range("N10").Value = 5.648181674564455
$ws.Range("B11").Value = 2.889857923294983
$ws.Range("C11").Value = 0.1175986038835504
$ws.Range("D11").Value = 0.1087364794088472
$ws.Range("E11").Value = 0.04752296404874912
$ws.Range("G11").Value = 0.002626154393822372
$ws.Range("K11").Value = 2.41765846226707
$ws.Range("L11").Value = 0.3079098453166296
$ws.Range("N11").Value = 5.718655742241538
$ws.Range("B12").Value = 2.913632502236794
$ws.Range("C12").Value = 0.1202398348067959
$ws.Range("D12").Value = 0.1109326642924913
$ws.Range("E12").Value = 0.04763616568208739
$ws.Range("G12").Value = 0.00262465777576276
$ws.Range("K12").Value = 2.441938652583644
$ws.Range("L12").Value = 0.3105084972751513
$ws.Range("N12").Value = 5.74541308086458
$ws.Range("B13").Value = 2.908501321283381
$ws.Range("C13").Value = 0.1196706540219452
$ws.Range("D13").Value = 0.1104593241099252
$ws.Range("E13").Value = 0.04761164581772803
$ws.Range("G13").Value = 0.002624978872860145
$ws.Range("K13").Value = 2.436700114078747
$ws.Range("L13").Value = 0.3099476555983784
$ws.Range("N13").Value = 5.739647256047931
$ws.Range("B14").Value = 2.891808963290828
$ws.Range("C14").Value = 0.117815744748782
$ws.Range("D14").Value = 0.1089170023225279
$ws.Range("E14").Value = 0.04753221419423781
$ws.Range("G14").Value = 0.002626030713536706
$ws.Range("K14").Value = 2.419651791136118
$ws.Range("L14").Value = 0.308123108601464
$ws.Range("N14").Value = 5.720855665974568
$ws.Range("B15").Value = 2.881616295405593
$ws.Range("C15").Value = 0.116680562895624
$ws.Range("D15").Value = 0.1079733143451307
$ws.Range("E15").Value = 0.04748396951061196
$ws.Range("G15").Value = 0.00262667858909848
$ws.Range("K15").Value = 2.409236586966017
$ws.Range("L15").Value = 0.3070089590901546
$ws.Range("N15").Value = 5.709354480870502
$ws.Range("B16").Value = 2.823693468818306
$ws.Range("C16").Value = 0.1101904606373694
$ws.Range("D16").Value = 0.1025810357740227
$ws.Range("E16").Value = 0.04721377188450582
$ws.Range("G16").Value = 0.002630446497774731
$ws.Range("K16").Value = 2.349969425143456
$ws.Range("L16").Value = 0.3006767525334624
$ws.Range("N16").Value = 5.643585786683616
$ws.Range("B17").Value = 2.788592891987491
$ws.Range("C17").Value = 0.1062222915081463
$ws.Range("D17").Value = 0.09928681885655521
$ws.Range("E17").Value = 0.04705360882100074
$ws.Range("G17").Value = 0.002632807371833162
$ws.Range("K17").Value = 2.313982362269314
$ws.Range("L17").Value = 0.2968388485424072
$ws.Range("N17").Value = 5.603362444000425
$ws.Range("B18").Value = 2.768562553925278
$ws.Range("C18").Value = 0.1039447002810334
$ws.Range("D18").Value = 0.09739706916380442
$ws.Range("E18").Value = 0.04696354196110875
$ws.Range("G18").Value = 0.00263418348543083
$ws.Range("K18").Value = 2.293419378300541
$ws.Range("L18").Value = 0.2946484799700073
$ws.Range("N18").Value = 5.580271975283608
$ws.Range("B19").Value = 2.761807805044896
$ws.Range("C19").Value = 0.1031743618076177
$ws.Range("D19").Value = 0.09675808453519608
$ws.Range("E19").Value = 0.04693339965272081
$ws.Range("G19").Value = 0.002634652544893039
$ws.Range("K19").Value = 2.286480369606011
$ws.Range("L19").Value = 0.293909787876018
$ws.Range("N19").Value = 5.572461609922101
$ws.Range("B20").Value = 2.792312981569808
$ws.Range("C20").Value = 0.1066442115785549
$ws.Range("D20").Value = 0.09963697528668547
$ws.Range("E20").Value = 0.04707044577207142
$ws.Range("G20").Value = 0.002632554170158121
$ws.Range("K20").Value = 2.317799174319759
$ws.Range("L20").Value = 0.2972456294262855
$ws.Range("N20").Value = 5.607639625833997
$ws.Range("B21").Value = 2.896705265498213
$ws.Range("C21").Value = 0.1183603668876856
$ws.Range("D21").Value = 0.1093698048863416
$ws.Range("E21").Value = 0.04755545984490972
$ws.Range("G21").Value = 0.002625721013931975
$ws.Range("K21").Value = 2.424653585823933
$ws.Range("L21").Value = 0.3086583053671887
$ws.Range("N21").Value = 5.726373290091999
$ws.Range("B22").Value = 2.966357018756469
$ws.Range("C22").Value = 0.1260622225173904
$ws.Range("D22").Value = 0.1157766383560812
$ws.Range("E22").Value = 0.04789076850931906
$ws.Range("G22").Value = 0.002621416112604802
$ws.Range("K22").Value = 2.495713027773434
$ws.Range("L22").Value = 0.316270846395895
$ws.Range("N22").Value = 5.804383322703472
$ws.Range("B23").Value = 2.929051498752187
$ws.Range("C23").Value = 0.1219474133029053
$ws.Range("D23").Value = 0.112352926907036
$ws.Range("E23").Value = 0.04771012994365442
$ws.Range("G23").Value = 0.002623699045465951
$ws.Range("K23").Value = 2.457674570833433
$ws.Range("L23").Value = 0.3121937533296517
$ws.Range("N23").Value = 5.762709803048324
$ws.Range("B24").Value = 2.790630663707304
$ws.Range("C24").Value = 0.1064534498187868
$ws.Range("D24").Value = 0.09947865668638656
$ws.Range("E24").Value = 0.04706282751670265
$ws.Range("G24").Value = 0.002632668584076292
$ws.Range("K24").Value = 2.316073199862103
$ws.Range("L24").Value = 0.2970616736772485
$ws.Range("N24").Value = 5.605705804687204
$ws.Range("B25").Value = 2.647376382923369
$ws.Range("C25").Value = 0.08992281830484217
$ws.Range("D25").Value = 0.08578186883910632
$ws.Range("E25").Value = 0.04644314528824545
$ws.Range("G25").Value = 0.002643040245885908
$ws.Range("K25").Value = 2.168517193931166
$ws.Range("L25").Value = 0.2813918379020635
$ws.Range("N25").Value = 5.438062713733103
